# Auto-generated edit script applying the crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.121.01"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.833.81"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6281"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07487"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.21"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07684"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "1.836.34"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.009"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6674"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.72"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009359"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -8.53%  "
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "29.114.27"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "2.082.22"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "223.14"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.100"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.01"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1390"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.487"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.499"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05727"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +9.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.148"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.081"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.205"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7407"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.831"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01779"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "1.213.52"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.516"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8891"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.05"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "1.981.40"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.51"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5096"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("B49").Value = "XinFinNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07641"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +12.81%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4064"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.974"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.40%  "
